$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column S that mirrors column R's formatting (this is how the
# sheet already built columns D..R for each successive year), then fill
# in the 2022 figures.
$ws.Range("R:R").Copy() | Out-Null
$ws.Range("S:S").Insert() | Out-Null

# Header year
$ws.Range("S4").Value = 2022

# Republic-level totals
$ws.Range("S5").Value = 135
$ws.Range("S6").Value = 99
$ws.Range("S7").Value = 36

# Batken oblast
$ws.Range("S8").Value = 97
$ws.Range("S9").Value = 80
$ws.Range("S10").Value = 17

# Jalal-Abad oblast
$ws.Range("S11").Value = 17
$ws.Range("S12").Value = 11
$ws.Range("S13").Value = 6

# Ysyk-Kul oblast
$ws.Range("S14").Value = 5
$ws.Range("S15").Value = 3
$ws.Range("S16").Value = 2

# Naryn oblast (rows 17-19 stay "-", already copied from column R)

# Osh oblast
$ws.Range("S20").Value = 6
$ws.Range("S21").Value = 1
$ws.Range("S22").Value = 5

# Talas oblast -> no incidents, becomes "-"
$ws.Range("S23").Value = "-"
$ws.Range("S24").Value = "-"
$ws.Range("S25").Value = "-"

# Chui oblast
$ws.Range("S26").Value = 10
$ws.Range("S27").Value = 4
$ws.Range("S28").Value = 6

# Bishkek city, Osh city (rows 29-34 stay "-", already copied from column R)

# Excel drops the cursor on the cell after the newly inserted/edited range
$ws.Range("T4").Select() | Out-Null
